$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Prediction" column inserted before the existing "Actual" column,
# shifting "Actual" from F to G, plus a new "Residual" column at H.
$ws.Range("F1").Value = "Prediction"
$ws.Range("G1").Value = "Actual"
$ws.Range("H1").Value = "Residual"

# Data row
$ws.Range("F2").Value = 6.7874299999999996
$ws.Range("H2").Formula = "=G2-F2"

# Column F width (new column) - engine quantizes to 1/6 character-width
# steps, so 13.3 is the closest achievable input to the recorded 14.1796875.
$ws.Columns("F").ColumnWidth = 13.3

# Match the recorded selection/cursor position after the edit.
$ws.Range("H3").Select()
